$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers to use the "<formatversion>" suffix instead of
#    the old "_old" / "_new" suffixes (e.g. "Segmentname_old" ->
#    "Segmentname_FV2210", "Segmentname_new" -> "Segmentname_FV2304").
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2210"
}
# Column 11 (K) is "diff" and keeps its name.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the data range into an Excel Table ("Table1") so the headers get
#    the autofilter dropdowns / structured reference support.
#    The header row already carries the bold/shaded/bordered style (s="1"),
#    so strip the explicit formatting before creating the table (otherwise
#    the engine records a headerRowDxfId) and then restore the same look
#    manually afterwards.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $dataRange, [Type]::Missing, 1)
$tbl.TableStyle = [Type]::Missing

$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
